# Reorder the "Recorded By" (column G) list of names/emails for every data
# row: move the last comma-separated entry to the front of the list.
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System")
# Rows whose value has only a single entry (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val.Split(",")
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
    $newVal = $rotated -join ", "

    $cell.Value2 = $newVal
}
